$wb = $excel.ActiveWorkbook

# Update zh-cn sheet's "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the a4d49d95-... row (row 4) to reflect a newly generated handback report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-22 20:47:37"
$wsZhCn.Range("H4").Value = "2016-03-22 20:48:04"

# Update de-de sheet's corresponding datetimes for the same row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-22 20:47:41"
$wsDeDe.Range("H4").Value = "2016-03-22 20:48:15"
